# Applies the "Added 3 APIs-Heartbeat, logoff, session" edit to worklistInfo.xlsx
# Changes (per the worklistCalibratorsControls sheet):
#   A3  : dev_worklist1    -> last1
#   A5  : Assert400        -> Assert404
#   A7  : dev_worklist1000 -> last1000001
#   A11 : dev_worklist1    -> last1
# Also moves the active cell selection from A11 to A14 to match the
# saved sheet view state in the target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("worklistCalibratorsControls")

$ws.Range("A3").Value = "last1"
$ws.Range("A5").Value = "Assert404"
$ws.Range("A7").Value = "last1000001"
$ws.Range("A11").Value = "last1"

$ws.Activate()
$ws.Range("A14").Select()
